$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (kemlabels-bulk-order-template) ---

# B1: service speed value text changed from "USPS Ground Advantage" to new descriptive label
$ws1.Range("B1").Value = "Ground Advantage: 1-5 days"

# Column B width changed (no longer best-fit, explicit custom width)
$ws1.Columns.Item(2).ColumnWidth = 27.83

# Data validation formula for courier list fixed to match named ranges (no spaces)
$ws1.Range("A1").Validation.Modify(3, 1, 1, '"USPS, UPSUSA, UPSCA"')

# Active selection on sheet1 moves to C14
$ws1.Activate()
$ws1.Range("C14").Select() | Out-Null

# --- Sheet2 ("Service Speeds") ---
$ws2.Range("A1").Value = "Ground Advantage: 1-5 days"
$ws2.Range("A2").Value = "Priority: 1-3 days"
$ws2.Range("A3").Value = "Express: 1-2 days"

$ws2.Range("C1").Value = "Express Early: 1 day"

$ws2.Range("B1").Value = "Next Day Air Early: 1 day"
$ws2.Range("B2").Value = "Next Day Air: 1 day"
$ws2.Range("B3").Value = "2nd Day Air: 2 days"
$ws2.Range("B4").Value = "3 Day Select: 3 days"
$ws2.Range("B5").Value = "Ground: Min 3 days"

$ws2.Range("C2").Value = "Express: 1 day"
$ws2.Range("C3").Value = "Express Saver: 1 day"
$ws2.Range("C4").Value = "Expedited: 2 days"
$ws2.Range("C5").Value = "Standard: Flexible"

# Active selection on sheet2 moves to C7
$ws2.Activate()
$ws2.Range("C7").Select() | Out-Null

# Re-activate sheet1 as the selected tab
$ws1.Activate()
